$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 246, pushing existing rows 246:259 down to 247:260
$ws.Rows("246:246").Insert()

# Populate the new row 246 with the new record's data.
# (Columns A,B,C,E,F,G,H,I,N,O,Q,R match the data already present in the row that was
#  pushed down to 247, since this is the same market/product series; only the
#  date / volume / price / Precio-$-Kg columns differ for the new record.)
$ws.Range("A246").Value2 = 11
$ws.Range("B246").Value2 = "Vega Monumental Concepción"
$ws.Range("C246").Value2 = "Bíobío"
$ws.Range("D246").Value2 = 45267
$ws.Range("D246").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E246").Value2 = 8
$ws.Range("F246").Value2 = 100112043
$ws.Range("G246").Value2 = "Pepino ensalada"
$ws.Range("H246").Value2 = "Sin especificar"
$ws.Range("I246").Value2 = "Primera"
$ws.Range("J246").Value2 = 100
$ws.Range("K246").Value2 = 15000
$ws.Range("L246").Value2 = 16000
$ws.Range("M246").Value2 = 15500
$ws.Range("N246").Value2 = "$/caja 60 unidades"
$ws.Range("O246").Value2 = "Región de Arica y Parinacota"
$ws.Range("P246").Value2 = 258
$ws.Range("Q246").Value2 = 60
$ws.Range("R246").Value2 = "Hortaliza"
